# ===================================================================
# Edit: add "Player Info" sheet, rename MATCH_CARD_LINK -> MATCH_CODE
# (storing bare match codes instead of full URLs) on "ODI Batting" and
# "ODI Bowling", drop the now-redundant blank INNING_NUMBER cells on
# "ODI Batting", and append a new "ODI Batting Extra" sheet.
# ===================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1. New first sheet: "Player Info"
#    (sheet handles in this host are positional, so re-fetch by name
#    right before use instead of caching references across structural
#    changes such as Add()/rename)
# -------------------------------------------------------------------
$battingForInsert = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingForInsert)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value2 = "ID"
$playerInfo.Range("B1").Value2 = "NAME"
$playerInfo.Range("C1").Value2 = "BATTING_HAND"
$playerInfo.Range("D1").Value2 = "BOWL_STYLE"

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value2 = "4062"
$playerInfo.Range("B2").Value2 = "Jasprit Jasbirsingh Bumrah"
$playerInfo.Range("C2").Value2 = "Right Handed"
$playerInfo.Range("D2").Value2 = "Right Arm Fast"

# -------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK (col D) -> MATCH_CODE, URL -> code,
#    drop the empty INNING_NUMBER (col B) placeholder cells.
# -------------------------------------------------------------------
$battingRows = $batting.UsedRange.Rows.Count

$batting.Range("D1").Value2 = "MATCH_CODE"

for ($r = 2; $r -le $battingRows; $r++) {
    $link = $batting.Cells.Item($r, 4).Value2
    $code = $link -replace ".*MatchCode=", ""
    $batting.Cells.Item($r, 4).NumberFormat = "@"
    $batting.Cells.Item($r, 4).Value2 = $code

    $inning = $batting.Cells.Item($r, 2).Value2
    if ($inning -eq "" -or $inning -eq $null) {
        $batting.Cells.Item($r, 2).ClearContents()
    }
}

# -------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK (col B) -> MATCH_CODE, URL -> code.
# -------------------------------------------------------------------
$bowlingRows = $bowling.UsedRange.Rows.Count

$bowling.Range("B1").Value2 = "MATCH_CODE"

for ($r = 2; $r -le $bowlingRows; $r++) {
    $link = $bowling.Cells.Item($r, 2).Value2
    $code = $link -replace ".*MatchCode=", ""
    $bowling.Cells.Item($r, 2).NumberFormat = "@"
    $bowling.Cells.Item($r, 2).Value2 = $code
}

# -------------------------------------------------------------------
# 4. New last sheet: "ODI Batting Extra"
# -------------------------------------------------------------------
$extra = $wb.Worksheets.Add($null, $bowling)
$extra.Name = "ODI Batting Extra"

$extra.Range("A1").Value2 = "MATCH_CODE"
$extra.Range("B1").Value2 = "BATTING_POSITION"
$extra.Range("C1").Value2 = "NUM_4"
$extra.Range("D1").Value2 = "NUM_6"
$extra.Range("E1").Value2 = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value2 = "MAN_OF_MATCH"

$extra.Columns.Item(1).NumberFormat = "@"
$extra.Columns.Item(3).NumberFormat = "@"
$extra.Columns.Item(4).NumberFormat = "@"
$extra.Columns.Item(5).NumberFormat = "@"
$extra.Columns.Item(6).NumberFormat = "@"

$extraRows = @(
    @("4332", $null, $null, $null, $null, "NO"),
    @("4338", 11, $null, $null, $null, "NO"),
    @("4342", 10, "0", "0", $null, "NO"),
    @("4345", 10, $null, $null, $null, "NO"),
    @("4350", $null, $null, $null, $null, "NO"),
    @("4353", 11, "0", "0", $null, "NO"),
    @("4398", 11, $null, $null, $null, "NO"),
    @("4399", $null, $null, $null, $null, "NO"),
    @("4400", $null, $null, $null, $null, "NO"),
    @("4402", 11, $null, $null, $null, "NO"),
    @("4406", $null, $null, $null, $null, "NO"),
    @("4410", $null, $null, $null, $null, "NO"),
    @("4435", 10, $null, $null, $null, "NO"),
    @("4436", $null, $null, $null, $null, "NO"),
    @("4437", 10, $null, $null, $null, "NO"),
    @("4524", 9, "2", "0", "4.24%", "NO"),
    @("4526", 9, $null, $null, $null, "YES"),
    @("4529", 9, "0", "0", "1.37%", "NO"),
    @("4609", $null, $null, $null, $null, $null),
    @("4613", $null, $null, $null, $null, $null)
)

$rowNum = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($rowNum, 1).Value2 = $row[0]

    if ($row[1] -ne $null) {
        $extra.Cells.Item($rowNum, 2).Value2 = $row[1]
    }
    if ($row[2] -ne $null) {
        $extra.Cells.Item($rowNum, 3).Value2 = $row[2]
    }
    if ($row[3] -ne $null) {
        $extra.Cells.Item($rowNum, 4).Value2 = $row[3]
    }
    if ($row[4] -ne $null) {
        $extra.Cells.Item($rowNum, 5).Value2 = $row[4]
    }
    if ($row[5] -ne $null) {
        $extra.Cells.Item($rowNum, 6).Value2 = $row[5]
    }
    $rowNum++
}
